$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Cardano -> USDC
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'

# Row 10: USDC -> Cardano
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'

# Row 2
$ws.Range('D2').Value = '96.363.23'
$ws.Range('E2').Value = '  +1.38%  '

# Row 3
$ws.Range('D3').Value = '3.573.24'
$ws.Range('E3').Value = '  -0.83%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.68'
$ws.Range('E5').Value = '  +1.96%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '653.42'
$ws.Range('E6').Value = '  -0.70%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.62'
$ws.Range('E7').Value = '  +11.74%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.409'
$ws.Range('E8').Value = '  +1.85%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.00%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.05'
$ws.Range('E10').Value = '  +6.84%  '

# Row 11
$ws.Range('D11').Value = '3.571.65'
$ws.Range('E11').Value = '  -0.78%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.33'
$ws.Range('E12').Value = '  +2.71%  '

# Row 13
$ws.Range('E13').Value = '  +1.39%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.38'
$ws.Range('E14').Value = '  +1.74%  '

# Row 15
$ws.Range('D15').Value = '4.240.64'
$ws.Range('E15').Value = '  -1.18%  '

# Row 16
$ws.Range('D16').Value = '96.233.89'
$ws.Range('E16').Value = '  +1.12%  '

# Row 17
$ws.Range('E17').Value = '  +2.93%  '

# Row 18
$ws.Range('D18').Value = '3.568.56'
$ws.Range('E18').Value = '  -1.30%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.75'
$ws.Range('E19').Value = '  -1.73%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.68'
$ws.Range('E20').Value = '  -1.31%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.80'
$ws.Range('E21').Value = '  -0.75%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.514'
$ws.Range('E22').Value = '  +7.75%  '

# Row 23
$ws.Range('E23').Value = '  -4.32%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '504.46'
$ws.Range('E24').Value = '  +0.28%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.94'
$ws.Range('E25').Value = '  +5.87%  '

# Row 26
$ws.Range('E26').Value = '  +2.44%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.14'
$ws.Range('E27').Value = '  +0.72%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.08'
$ws.Range('E28').Value = '  +4.98%  '

# Row 29
$ws.Range('D29').Value = '3.765.95'
$ws.Range('E29').Value = '  -0.87%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.155'
$ws.Range('E30').Value = '  +12.91%  '

# Row 31
$ws.Range('E31').Value = '  -4.35%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.41'
$ws.Range('E32').Value = '  +1.59%  '

# Row 33
$ws.Range('E33').Value = '  -0.06%  '

# Row 34
$ws.Range('E34').Value = '  +3.84%  '

# Row 35
$ws.Range('E35').Value = '  +0.29%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.35'
$ws.Range('E36').Value = '  -1.22%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '635.27'
$ws.Range('E37').Value = '  +10.70%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.87'
$ws.Range('E38').Value = '  +10.33%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.565'
$ws.Range('E39').Value = '  +1.92%  '

# Row 40
$ws.Range('E40').Value = '  +12.37%  '

# Row 42
$ws.Range('E42').Value = '  +1.06%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.904'
$ws.Range('E43').Value = '  -1.26%  '

# Row 44
$ws.Range('E44').Value = '  +6.25%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.72'
$ws.Range('E45').Value = '  +1.88%  '

# Row 46
$ws.Range('E46').Value = '  +3.82%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.29'
$ws.Range('E47').Value = '  +3.57%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.53'
$ws.Range('E48').Value = '  -0.53%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.19'
$ws.Range('E49').Value = '  -4.89%  '

# Row 50
$ws.Range('E50').Value = '  +0.83%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.23'
$ws.Range('E51').Value = '  +3.72%  '
